# ---------------------------------------------------------------------------
# Rename the "Requested quantity" header on the two existing sheets, and add
# a new "PO Forecast" sheet with a Prophet-style forecast table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet: B1 header rename -----------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet: B1 header rename --------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" sheet as the last (3rd) sheet -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$data = @(
    @(44983.99999999999, 282, 30.89946159354966, 539.7493863426918),
    @(44990.99999999999, 266, -13.70563179181091, 517.7505694592689),
    @(45011.99999999999, 218, -27.64063484482162, 472.3502427026586),
    @(45046.99999999999, 137, -106.9239103209763, 391.9131011521557),
    @(45081.99999999999, 57, -195.2510265283912, 317.0390309989954),
    @(45088.99999999999, 41, -203.965720741007, 301.0109286838519),
    @(45095.99999999999, 25, -232.6256260778344, 265.6382141296934),
    @(45102.99999999999, 8, -248.4313573036023, 272.2085128700767),
    @(45109.99999999999, 0, -262.5062798242312, 243.0419260745167),
    @(45116.99999999999, 0, -278.7548464782893, 248.7365052565403),
    @(45123.99999999999, 0, -309.6129174705306, 212.0131762967741),
    @(45130.99999999999, 0, -312.4929977414869, 204.2336894513768),
    @(45137.99999999999, 0, -314.746198841325, 180.6242076601472),
    @(45144.99999999999, 0, -351.6790070733148, 166.4423996117111)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Range("A$r").Value = $row[0]
    $wsForecast.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Range("B$r").Value = $row[1]
    $wsForecast.Range("C$r").Value = $row[2]
    $wsForecast.Range("D$r").Value = $row[3]
    $r++
}

